$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 110, shifting existing rows 110-162 down to 111-163
$ws.Rows("110:110").Insert()

# Fill in the new row 110 with the new data (same structure as surrounding rows)
$ws.Cells.Item(110, 1).Value = 8
$ws.Cells.Item(110, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(110, 3).Value = "Coquimbo"
$ws.Cells.Item(110, 4).Value = 45233
$ws.Cells.Item(110, 5).Value = 4
$ws.Cells.Item(110, 6).Value = 100112028
$ws.Cells.Item(110, 7).Value = "Sandia"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 1200
$ws.Cells.Item(110, 11).Value = 650
$ws.Cells.Item(110, 12).Value = 700
$ws.Cells.Item(110, 13).Value = 675
$ws.Cells.Item(110, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(110, 15).Value = "Perú"
$ws.Cells.Item(110, 16).Value = 675
$ws.Cells.Item(110, 17).Value = 1
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Match the date-style formatting of column D used elsewhere (numFmt 165)
$ws.Cells.Item(110, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
